$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column (H): copy formatting from the neighboring header cell (G1)
# so the new header cell matches the existing header style, then set its value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Populate the new column's data values for the two data rows
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
